$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 127
$ws.Range("B3").Value = 125
$ws.Range("B38").Value = 133
$ws.Range("B48").Value = 132

$ws.Range("F46").Select()
